$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPaths = @(
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.38.32.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.11.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.18.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.25.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.29.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.34.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.49.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.52.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.39.57.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.17.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.21.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.24.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.27.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.30.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.33.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.37.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.55.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.40.58.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.41.01.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.41.04.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.41.20.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.41.23.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.41.25.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.41.27.png',
    '/Users/KevinGao/Desktop/fiber_clot/fiber/16x16_modified/trainingtraining_set2.41.30.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.50.55.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.01.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.07.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.11.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.15.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.19.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.23.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.26.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.30.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.34.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.37.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.41.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.44.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.51.48.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.26.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.30.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.34.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.37.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.40.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.43.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.48.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.50.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.53.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.55.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.57.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.52.59.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.53.02.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.53.04.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.53.08.png',
    '/Users/KevinGao/Desktop/background/16x16_modified16x16/trainingtraining_set9.53.14.png'
)

for ($i = 0; $i -lt $newPaths.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = $newPaths[$i]
}

# Update the sheet view: keep gridlines on (matches original workbook setting),
# scroll the top-left visible cell to A17, and move the active selection to C23.
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("C23").Select() | Out-Null

